$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking values (e.g. "320.19")
# are stored as strings, matching the source data which used inline strings.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '48.945.82'
$ws.Range("E2").Value = '  +1.63%  '
$ws.Range("D3").Value = '2.530.27'
$ws.Range("E3").Value = '  +0.88%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = '320.19'
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").Value = '107.07'
$ws.Range("E6").Value = '  -1.79%  '
$ws.Range("D7").Value = '0.524'
$ws.Range("E7").Value = '  -0.96%  '
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").Value = '0.546'
$ws.Range("E9").Value = '  +0.14%  '
$ws.Range("D10").Value = '39.58'
$ws.Range("E10").Value = '  -0.68%  '
$ws.Range("D11").Value = '20.14'
$ws.Range("E11").Value = '  +0.13%  '
$ws.Range("D12").Value = '0.0808'
$ws.Range("E12").Value = '  -1.23%  '
$ws.Range("D13").Value = '0.126'
$ws.Range("E13").Value = '  +0.99%  '
$ws.Range("D14").Value = '7.16'
$ws.Range("E14").Value = '  -0.48%  '
$ws.Range("D15").Value = '2.926.10'
$ws.Range("E15").Value = '  +0.76%  '
$ws.Range("D16").Value = '2.524.09'
$ws.Range("E16").Value = '  +0.54%  '
$ws.Range("D17").Value = '0.849'
$ws.Range("E17").Value = '  +0.29%  '
$ws.Range("D18").Value = '48.756.29'
$ws.Range("E18").Value = '  +1.55%  '
$ws.Range("D19").Value = '12.98'
$ws.Range("E19").Value = '  -1.83%  '
$ws.Range("D20").Value = '2.93'
$ws.Range("E20").Value = '  +7.77%  '
$ws.Range("D21").Value = '6.60'
$ws.Range("E21").Value = '  -0.15%  '
$ws.Range("D22").Value = '0.0₃0937'
$ws.Range("E22").Value = '  -0.51%  '
$ws.Range("D23").Value = '286.66'
$ws.Range("E23").Value = '  +4.41%  '
$ws.Range("D24").Value = '71.31'
$ws.Range("E24").Value = '  -1.25%  '
$ws.Range("D25").Value = '2.50'
$ws.Range("E25").Value = '  -2.42%  '
$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").Value = '26.06'
$ws.Range("E26").Value = '  +0.65%  '
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.15%  '
$ws.Range("D28").Value = '2.23'
$ws.Range("E28").Value = '  -7.23%  '
$ws.Range("D29").Value = '0.144'
$ws.Range("E29").Value = '  +2.36%  '
$ws.Range("D30").Value = '9.72'
$ws.Range("E30").Value = '  -3.51%  '
$ws.Range("D31").Value = '35.15'
$ws.Range("E31").Value = '  -1.00%  '
$ws.Range("D32").Value = '49.34'
$ws.Range("E32").Value = '  -0.21%  '
$ws.Range("D33").Value = '19.50'
$ws.Range("E33").Value = '  +0.82%  '
$ws.Range("E34").Value = '  -0.33%  '
$ws.Range("D35").Value = '5.31'
$ws.Range("E35").Value = '  -0.71%  '
$ws.Range("D36").Value = '0.0775'
$ws.Range("E36").Value = '  -1.04%  '
$ws.Range("D37").Value = '1.98'
$ws.Range("E37").Value = '  +0.97%  '
$ws.Range("D38").Value = '4.59'
$ws.Range("E38").Value = '  -0.68%  '
$ws.Range("E39").Value = '  -1.27%  '
$ws.Range("E40").Value = '  -0.66%  '
$ws.Range("E41").Value = '  +0.17%  '
$ws.Range("D42").Value = '120.34'
$ws.Range("E42").Value = '  -1.92%  '
$ws.Range("D43").Value = '21.84'
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("D44").Value = '0.0304'
$ws.Range("E44").Value = '  +0.09%  '
$ws.Range("D45").Value = '2.006.69'
$ws.Range("E45").Value = '  -0.84%  '
$ws.Range("D46").Value = '3.19'
$ws.Range("E46").Value = '  +1.83%  '
$ws.Range("D47").Value = '1.99'
$ws.Range("E47").Value = '  +7.65%  '
$ws.Range("E48").Value = '  +6.21%  '
$ws.Range("E49").Value = '  -0.18%  '
$ws.Range("D50").Value = '5.20'
$ws.Range("E50").Value = '  +0.24%  '
$ws.Range("D51").Value = '80.59'
$ws.Range("E51").Value = '  +1.49%  '

# Restore default cell style on column D (NumberFormat change above
# would otherwise leave a stray "Text" style on these cells).
$dRange.Style = "Normal"

